{"js": "// Replace the date line and each arithmetic problem's text with its\n// updated value. Every <w:t> run in this document is a unique string,\n// so a simple search-and-replace per value is safe and unambiguous.\nconst replacements = [\n  [\"2025-12-31 Wednesday\", \"2026-01-01 Thursday\"],\n  [\"624\u00f74=\", \"768\u00f75=\"],\n  [\"268\u00f79=\", \"690\u00f79=\"],\n  [\"681\u00f74=\", \"350\u00f79=\"],\n  [\"112\u00f74=\", \"240\u00f77=\"],\n  [\"702\u00f76=\", \"355\u00f77=\"],\n  [\"516\u00f79=\", \"527\u00f72=\"],\n  [\"229\u00f78=\", \"268\u00f75=\"],\n  [\"665\u00f77=\", \"274\u00f77=\"],\n  [\"467\u00f77=\", \"777\u00f73=\"],\n  [\"537\u00f74=\", \"116\u00f72=\"],\n  [\"226\u00f75=\", \"874\u00f78=\"],\n  [\"875\u00f77=\", \"188\u00f75=\"],\n  [\"483\u00f72=\", \"438\u00f74=\"],\n  [\"578\u00f72=\", \"793\u00f79=\"],\n  [\"947\u00f79=\", \"660\u00f76=\"],\n  [\"223\u00f72=\", \"850\u00f78=\"],\n  [\"794\u00f77=\", \"821\u00f74=\"],\n  [\"163\u00f79=\", \"968\u00f74=\"],\n  [\"725\u00f76=\", \"880\u00f75=\"],\n  [\"319\u00f77=\", \"659\u00f78=\"],\n  [\"805\u00f79=\", \"589\u00f77=\"],\n  [\"421\u00f77=\", \"195\u00f77=\"],\n  [\"283\u00f74=\", \"956\u00f76=\"],\n  [\"439\u00f73=\", \"843\u00f78=\"],\n  [\"798\u00f74=\", \"848\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each arithmetic problem's text with its\n# updated value. Every piece of text in this document is unique, so a\n# simple Find/Replace pass per value is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"2025-12-31 Wednesday\"; new=\"2026-01-01 Thursday\"},\n    @{old=\"624\u00f74=\"; new=\"768\u00f75=\"},\n    @{old=\"268\u00f79=\"; new=\"690\u00f79=\"},\n    @{old=\"681\u00f74=\"; new=\"350\u00f79=\"},\n    @{old=\"112\u00f74=\"; new=\"240\u00f77=\"},\n    @{old=\"702\u00f76=\"; new=\"355\u00f77=\"},\n    @{old=\"516\u00f79=\"; new=\"527\u00f72=\"},\n    @{old=\"229\u00f78=\"; new=\"268\u00f75=\"},\n    @{old=\"665\u00f77=\"; new=\"274\u00f77=\"},\n    @{old=\"467\u00f77=\"; new=\"777\u00f73=\"},\n    @{old=\"537\u00f74=\"; new=\"116\u00f72=\"},\n    @{old=\"226\u00f75=\"; new=\"874\u00f78=\"},\n    @{old=\"875\u00f77=\"; new=\"188\u00f75=\"},\n    @{old=\"483\u00f72=\"; new=\"438\u00f74=\"},\n    @{old=\"578\u00f72=\"; new=\"793\u00f79=\"},\n    @{old=\"947\u00f79=\"; new=\"660\u00f76=\"},\n    @{old=\"223\u00f72=\"; new=\"850\u00f78=\"},\n    @{old=\"794\u00f77=\"; new=\"821\u00f74=\"},\n    @{old=\"163\u00f79=\"; new=\"968\u00f74=\"},\n    @{old=\"725\u00f76=\"; new=\"880\u00f75=\"},\n    @{old=\"319\u00f77=\"; new=\"659\u00f78=\"},\n    @{old=\"805\u00f79=\"; new=\"589\u00f77=\"},\n    @{old=\"421\u00f77=\"; new=\"195\u00f77=\"},\n    @{old=\"283\u00f74=\"; new=\"956\u00f76=\"},\n    @{old=\"439\u00f73=\"; new=\"843\u00f78=\"},\n    @{old=\"798\u00f74=\"; new=\"848\u00f72=\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($r.old, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
